$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column D; existing D:K data shifts right to E:L.
$ws.Columns("D").Insert()

# Copy formatting (number formats/styles) from the shifted former-D column (now E) back onto new D,
# restricted to the populated row blocks so the sheet dimension/used range is not inflated.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the new (most recent) reporting-period figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 8611000
$ws.Range("D9").Value = 3339000
$ws.Range("D10").Value = 5272000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 675000
$ws.Range("D15").Value = 1954000
$ws.Range("D17").Value = 7783000
$ws.Range("D18").Value = 828000
$ws.Range("D20").Value = 3000
$ws.Range("D21").Value = 2785000
$ws.Range("D22").Value = 1536000
$ws.Range("D23").Value = -705000
$ws.Range("D24").Value = -62000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -643000
$ws.Range("D27").Value = -750000
$ws.Range("D28").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3000
$ws.Range("D33").Value = -750000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -750000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 354000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 783000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 193000
$ws.Range("D46").Value = 1330000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 14187000
$ws.Range("D49").Value = 7877000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 265000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 23659000
$ws.Range("D57").Value = 495000
$ws.Range("D58").Value = 844000
$ws.Range("D59").Value = 1222000
$ws.Range("D60").Value = 2561000
$ws.Range("D61").Value = 16427000
$ws.Range("D62").Value = 3071000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 22059000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -2752000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1600000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -750000
$ws.Range("D83").Value = 1954000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1812000
$ws.Range("D91").Value = -1192000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1176000
$ws.Range("D96").Value = -107000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -608000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 28000

# These rows carry the "NA" marker in column D instead of a numeric figure.
$ws.Range("D12").Value = "NA"
$ws.Range("D29").Value = "NA"

# A handful of rows also had their (now shifted) column E figure restated/corrected.
$ws.Range("E9").Value = 3555000
$ws.Range("E10").Value = 5573000
$ws.Range("E14").Value = 3028000
$ws.Range("E20").Value = 1537000
$ws.Range("E21").Value = 2065000
$ws.Range("E22").Value = 3068000
$ws.Range("E32").Value = -1537000

Write-Output "done"
